# Update cryptos list cell values (price & volume columns) per source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.711.48"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.722.29"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'239.77"
$ws.Range("D6").Value = "'0.9988"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.4830"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").Value = "'0.2577"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "1.724.94"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "'0.06872"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'0.6037"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'4.460"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "'76.83"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "26.538.35"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'0.000007137"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "1.947.07"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'4.413"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "'8.570"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "'5.053"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'139.09"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'15.23"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'1.769"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").Value = "'106.30"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "'4.020"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'0.07906"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Value = "'3.664"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'0.04474"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "'0.9978"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'2.596"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'0.9989"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'0.6166"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").Value = "'0.9257"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'2.006"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").Value = "'2.441"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").Value = "'0.9983"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.01490"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "'5.610"
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("D44").Value = "'99.90"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'0.3817"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'6.773"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'0.1150"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "'0.05369"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'7.836"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'30.02"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'1.233"
$ws.Range("E51").Value = "  +0.35%  "
